$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 2 and row 5 for columns D, L, M, N, O, P, S

# Row 2 (was: 44881, Segunda, 100, 11250, 11250, 11250, 11250)
# becomes: 44874, Primera, 200, 7500, 8000, 7750, 7750
$ws.Range("D2").Value = 44874
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 7500
$ws.Range("O2").Value = 8000
$ws.Range("P2").Value = 7750
$ws.Range("S2").Value = 7750

# Row 5 (was: 44874, Primera, 200, 7500, 8000, 7750, 7750)
# becomes: 44881, Segunda, 100, 11250, 11250, 11250, 11250
$ws.Range("D5").Value = 44881
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 11250
$ws.Range("O5").Value = 11250
$ws.Range("P5").Value = 11250
$ws.Range("S5").Value = 11250
